$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 42604.891412037039
$ws.Range("B5").Value = "Named"
$ws.Range("C5").Value = 13670
$ws.Range("D5").Value = 8177
$ws.Range("E5").Value = 470
$ws.Range("F5").Value = 74
$ws.Range("G5").Value = 65
$ws.Range("H5").Value = 52
$ws.Range("I5").Value = 46
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 6
$ws.Range("L5").Value = 14
$ws.Range("M5").Value = 85
